$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44545
$ws.Range("K2").Value = "Castle Brite"
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 19000
$ws.Range("P2").Value = 18500
$ws.Range("Q2").Value = "`$/caja 15 kilos"
$ws.Range("S2").Value = 1233
$ws.Range("T2").Value = 15

# Row 3
$ws.Range("D3").Value = 44545
$ws.Range("K3").Value = "Castle Brite"
$ws.Range("N3").Value = 17000
$ws.Range("O3").Value = 17000
$ws.Range("P3").Value = 17000
$ws.Range("Q3").Value = "`$/caja 15 kilos"
$ws.Range("S3").Value = 1133
$ws.Range("T3").Value = 15

# Row 4
$ws.Range("D4").Value = 44159
$ws.Range("K4").Value = "Castle Brite"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 14500
$ws.Range("Q4").Value = "`$/caja 15 kilos"
$ws.Range("S4").Value = 967
$ws.Range("T4").Value = 15

# Row 7
$ws.Range("D7").Value = 44559
$ws.Range("K7").Value = "Modesto"
$ws.Range("N7").Value = 19000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 19500
$ws.Range("Q7").Value = "`$/caja 18 kilos"
$ws.Range("S7").Value = 1083
$ws.Range("T7").Value = 18

# Row 8
$ws.Range("D8").Value = 44559
$ws.Range("K8").Value = "Modesto"
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 18000
$ws.Range("Q8").Value = "`$/caja 18 kilos"
$ws.Range("S8").Value = 1000
$ws.Range("T8").Value = 18

# Row 9
$ws.Range("D9").Value = 44187
$ws.Range("K9").Value = "Dina"
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 16000
$ws.Range("P9").Value = 15500
$ws.Range("Q9").Value = "`$/caja 18 kilos"
$ws.Range("S9").Value = 861
$ws.Range("T9").Value = 18

# Row 10
$ws.Range("D10").Value = 44579
$ws.Range("K10").Value = "Modesto"
$ws.Range("M10").Value = 180
$ws.Range("N10").Value = 13000
$ws.Range("O10").Value = 14000
$ws.Range("P10").Value = 13444
$ws.Range("S10").Value = 747
